# "Actualizado el plan de pruebas, todas las pruebas de la app OK"
#
# Test case 11 ("Modificar un vehiculo", row 12) had its status (column F)
# still marked as "PTE" (pending). Update it to "OK" and give it the same
# green "OK" formatting used by the other completed test cases (e.g. F9,
# which already carries the OK fill/border style) instead of the orange
# "PTE" formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an already-"OK" status cell onto F12, then set its
# value/text to "OK" so both the fill and the content reflect the update.
$ws.Range("F9").Copy()
$ws.Range("F12").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("F12").Value = "OK"

# Reflect the author's last active cell/selection in the sheet.
$ws.Range("G11").Select()
